# data_2023-11-24_LSR3_H.xlsx -- minor edits and comments
# - comments on methods checklist document (risk-of-bias "Overall"/D-domain
#   columns, rows 4-15)
# - chunks options added (dropout_any_e/_n in columns R/S for rows 6-11,
#   plus overall_baseline in column AH for rows 6-11)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: some randomized_n/response_n/dropout_any_n columns (N, O, R, S)
# store numbers as literal text (t="inlineStr") in this sheet, not real
# numbers -- force the cell to Text before assigning so the digits aren't
# auto-coerced into a numeric cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# ---- Row 4 (Koblan 2020, ulotaront) ----
$ws.Range("AE4").Value = "Low"

# ---- Row 5 (Koblan 2020, placebo) ----
$ws.Range("AE5").Value = "Low"

# ---- Row 6 (NCT04072354, ulotaront 50mg/d) ----
$ws.Range("J6").Value = 144
Set-TextValue $ws.Range("N6") "22"
Set-TextValue $ws.Range("O6") "144"
Set-TextValue $ws.Range("R6") "34"
Set-TextValue $ws.Range("S6") "144"
$ws.Range("AB6").Value = "Some concerns"
$ws.Range("AD6").Value = "Low"
$ws.Range("AE6").Value = "Low"
$ws.Range("AG6").Value = 142
$ws.Range("AH6").Value = 102
$ws.Range("AJ6").Value = 19.07
$ws.Range("AM6").Value = "yes"

# ---- Row 7 (NCT04072354, ulotaront 75mg/d) ----
Set-TextValue $ws.Range("N7") "29"
Set-TextValue $ws.Range("R7") "27"
Set-TextValue $ws.Range("S7") "145"
$ws.Range("AB7").Value = "Some concerns"
$ws.Range("AD7").Value = "Low"
$ws.Range("AE7").Value = "Low"
$ws.Range("AH7").Value = 102
$ws.Range("AJ7").Value = 19.27
$ws.Range("AM7").Value = "yes"

# ---- Row 8 (NCT04072354, placebo) ----
$ws.Range("J8").Value = 146
Set-TextValue $ws.Range("N8") "26"
Set-TextValue $ws.Range("O8") "146"
Set-TextValue $ws.Range("R8") "27"
Set-TextValue $ws.Range("S8") "146"
$ws.Range("AB8").Value = "Some concerns"
$ws.Range("AD8").Value = "Low"
$ws.Range("AE8").Value = "Low"
$ws.Range("AH8").Value = 102
$ws.Range("AJ8").Value = 18.06
$ws.Range("AM8").Value = "yes"

# ---- Row 9 (NCT04092686, ulotaront 75mg/d) ----
$ws.Range("J9").Value = 155
Set-TextValue $ws.Range("N9") "23"
Set-TextValue $ws.Range("O9") "155"
Set-TextValue $ws.Range("R9") "34"
Set-TextValue $ws.Range("S9") "155"
$ws.Range("AB9").Value = "Some concerns"
$ws.Range("AD9").Value = "Low"
$ws.Range("AE9").Value = "Low"
$ws.Range("AG9").Value = 153
$ws.Range("AH9").Value = 101
$ws.Range("AJ9").Value = 18.55
$ws.Range("AM9").Value = "yes"

# ---- Row 10 (NCT04092686, ulotaront 100mg/d) ----
Set-TextValue $ws.Range("N10") "27"
Set-TextValue $ws.Range("R10") "38"
Set-TextValue $ws.Range("S10") "154"
$ws.Range("AB10").Value = "Some concerns"
$ws.Range("AD10").Value = "Low"
$ws.Range("AE10").Value = "Low"
$ws.Range("AG10").Value = 152
$ws.Range("AH10").Value = 100
$ws.Range("AJ10").Value = 18.49
$ws.Range("AM10").Value = "yes"

# ---- Row 11 (NCT04092686, placebo) ----
$ws.Range("J11").Value = 155
Set-TextValue $ws.Range("N11") "21"
Set-TextValue $ws.Range("O11") "155"
Set-TextValue $ws.Range("R11") "27"
Set-TextValue $ws.Range("S11") "155"
$ws.Range("AB11").Value = "Some concerns"
$ws.Range("AD11").Value = "Low"
$ws.Range("AE11").Value = "Low"
$ws.Range("AG11").Value = 155
$ws.Range("AH11").Value = 100
$ws.Range("AJ11").Value = 18.67
$ws.Range("AM11").Value = "yes"

# ---- Row 12 (NCT04512066, placebo) ----
$ws.Range("AE12").Value = "Low"

# ---- Row 13 (NCT04512066, ralmitaront 45mg) ----
# death_e/serious_e (V13/X13) were transposed; 0 deaths, 2 serious events.
Set-TextValue $ws.Range("V13") "0"
Set-TextValue $ws.Range("X13") "2"
$ws.Range("AE13").Value = "Low"

# ---- Row 14 (NCT04512066, ralmitaront 90mg) ----
$ws.Range("AE14").Value = "Low"

# ---- Row 15 (NCT04512066, ralmitaront 180mg) ----
$ws.Range("AE15").Value = "Low"

# ---- Row 18 (Perini 2023, amisulpride) ----
# overall_baseline was an erroneous change score (-10.69), not a baseline
# -- remove it.
$ws.Range("AH18").ClearContents()

# ---- Row 19 (Tsukada 2023, placebo) ----
# same cleanup as row 18 (-17.3 was a change score, not a baseline).
$ws.Range("AH19").ClearContents()
